$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 16795.334
$ws.Range("I69").Value = 8777
$ws.Range("J69").Value = 18399
$ws.Range("K69").Value = 26331
$ws.Range("L69").Value = 55197
$ws.Range("M69").Value = -25457
$ws.Range("N69").Value = -56945
$ws.Range("H72").Value = 16795.334
$ws.Range("I72").Value = 8777
$ws.Range("J72").Value = 18399
$ws.Range("K72").Value = 78993
$ws.Range("L72").Value = 165591
$ws.Range("M72").Value = -74625
$ws.Range("N72").Value = -174327
$ws.Range("H86").Value = 170281.5
$ws.Range("I86").Value = 203897.8
$ws.Range("K86").Value = 203897.8
$ws.Range("M86").Value = -202774.8
$ws.Range("H89").Value = 170281.5
$ws.Range("I89").Value = 203897.8
$ws.Range("K89").Value = 1019489
$ws.Range("M89").Value = -1013873
$ws.Range("H96").Value = 2424.5386
$ws.Range("I96").Value = 579
$ws.Range("J96").Value = 3578
$ws.Range("K96").Value = 1737
$ws.Range("L96").Value = 10734
$ws.Range("M96").Value = -364
$ws.Range("N96").Value = -13480
$ws.Range("H98").Value = 1410.5
$ws.Range("I98").Value = 1146.4117
$ws.Range("J98").Value = 5900
$ws.Range("K98").Value = 1146.4117
$ws.Range("L98").Value = 5900
$ws.Range("M98").Value = 351.5882999999999
$ws.Range("N98").Value = -8896
$ws.Range("H115").Value = 1163.4546
$ws.Range("I115").Value = 279.33334
$ws.Range("K115").Value = 838.0000200000001
$ws.Range("M115").Value = 728.9999799999999
$ws.Range("H116").Value = 5748.8975
$ws.Range("I116").Value = 5231.483
$ws.Range("J116").Value = 7249.4
$ws.Range("K116").Value = 5231.483
$ws.Range("L116").Value = 7249.4
$ws.Range("M116").Value = -1789.483
$ws.Range("N116").Value = -14133.4
$ws.Range("H122").Value = 1410.5
$ws.Range("I122").Value = 1146.4117
$ws.Range("J122").Value = 5900
$ws.Range("K122").Value = 3439.2351
$ws.Range("L122").Value = 17700
$ws.Range("M122").Value = -989.2351000000003
$ws.Range("N122").Value = -22600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3529.3333
$ws.Range("I45").Value = 1088.8334
$ws.Range("J45").Value = 4505.533
$ws.Range("K45").Value = 1088.8334
$ws.Range("L45").Value = 4505.533
$ws.Range("M45").Value = -711.8334
$ws.Range("N45").Value = -5259.533
$ws.Range("H97").Value = 7134.696
$ws.Range("I97").Value = 9640.308
$ws.Range("J97").Value = 3877.4
$ws.Range("K97").Value = 9640.308
$ws.Range("L97").Value = 3877.4
$ws.Range("M97").Value = -9144.308
$ws.Range("N97").Value = -4869.4
$ws.Range("H122").Value = 1227.76
$ws.Range("I122").Value = 1070.5834
$ws.Range("K122").Value = 3211.7502
$ws.Range("M122").Value = -761.7502
$ws.Range("H132").Value = 38927.965
$ws.Range("I132").Value = 40083.348
$ws.Range("K132").Value = 120250.044
$ws.Range("M132").Value = -117720.044

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 34612.062
$ws.Range("I99").Value = 46075
$ws.Range("K99").Value = 46075
$ws.Range("M99").Value = -44577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4918.909
$ws.Range("I31").Value = 3333.4
$ws.Range("K31").Value = 3333.4
$ws.Range("M31").Value = -3038.4
$ws.Range("H34").Value = 4918.909
$ws.Range("I34").Value = 3333.4
$ws.Range("K34").Value = 3333.4
$ws.Range("M34").Value = -3131.4
$ws.Range("H58").Value = 69694
$ws.Range("I58").Value = 93965.18
$ws.Range("J58").Value = 2948.25
$ws.Range("K58").Value = 93965.18
$ws.Range("L58").Value = 2948.25
$ws.Range("M58").Value = -93762.18
$ws.Range("N58").Value = -3354.25
$ws.Range("H86").Value = 20787.738
$ws.Range("I86").Value = 39915.273
$ws.Range("J86").Value = 3254.1667
$ws.Range("K86").Value = 39915.273
$ws.Range("L86").Value = 3254.1667
$ws.Range("M86").Value = -38792.273
$ws.Range("N86").Value = -5500.1667
$ws.Range("H89").Value = 20787.738
$ws.Range("I89").Value = 39915.273
$ws.Range("J89").Value = 3254.1667
$ws.Range("K89").Value = 199576.365
$ws.Range("L89").Value = 16270.8335
$ws.Range("M89").Value = -193960.365
$ws.Range("N89").Value = -27502.8335
$ws.Range("H136").Value = 69694
$ws.Range("I136").Value = 93965.18
$ws.Range("J136").Value = 2948.25
$ws.Range("K136").Value = 281895.54
$ws.Range("L136").Value = 8844.75
$ws.Range("M136").Value = -279345.54
$ws.Range("N136").Value = -13944.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 2365.6667
$ws.Range("I51").Value = 2365.6667
$ws.Range("K51").Value = 7097.000100000001
$ws.Range("M51").Value = -6637.000100000001
$ws.Range("H68").Value = 427.375
$ws.Range("J68").Value = 407.2
$ws.Range("L68").Value = 1221.6
$ws.Range("N68").Value = -2843.6
$ws.Range("H71").Value = 427.375
$ws.Range("J71").Value = 407.2
$ws.Range("L71").Value = 3664.8
$ws.Range("N71").Value = -11776.8
$ws.Range("H92").Value = 554.46155
$ws.Range("I92").Value = 426
$ws.Range("J92").Value = 760
$ws.Range("K92").Value = 1278
$ws.Range("L92").Value = 2280
$ws.Range("M92").Value = -30
$ws.Range("N92").Value = -4776
$ws.Range("H134").Value = 2974.7144
$ws.Range("I134").Value = 2004.6
$ws.Range("J134").Value = 5400
$ws.Range("K134").Value = 6013.799999999999
$ws.Range("L134").Value = 16200
$ws.Range("M134").Value = -943.7999999999993
$ws.Range("N134").Value = -26340
$ws.Range("H140").Value = 1556.7059
$ws.Range("I140").Value = 1556.7059
$ws.Range("K140").Value = 4670.1177
$ws.Range("M140").Value = 509.8823000000002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 13000
$ws.Range("J2").Value = 13000
$ws.Range("L2").Value = 13000
$ws.Range("N2").Value = -13224
$ws.Range("H22").Value = 65185.25
$ws.Range("I22").Value = 143672
$ws.Range("J22").Value = 4140
$ws.Range("K22").Value = 143672
$ws.Range("L22").Value = 4140
$ws.Range("M22").Value = -143377
$ws.Range("N22").Value = -4730
$ws.Range("H27").Value = 65185.25
$ws.Range("I27").Value = 143672
$ws.Range("J27").Value = 4140
$ws.Range("K27").Value = 143672
$ws.Range("L27").Value = 4140
$ws.Range("M27").Value = -143565
$ws.Range("N27").Value = -4354
$ws.Range("H123").Value = 63990.668
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 63990.668
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 63990.668
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -73790.668
$ws.Range("H136").Value = 3171.282
$ws.Range("I136").Value = 2712.5557
$ws.Range("J136").Value = 3564.476
$ws.Range("K136").Value = 8137.6671
$ws.Range("L136").Value = 10693.428
$ws.Range("M136").Value = -5587.6671
$ws.Range("N136").Value = -15793.428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 655
$ws.Range("I122").Value = 655
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1965
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 485
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 52222
$ws.Range("J131").Value = 52222
$ws.Range("L131").Value = 52222
$ws.Range("N131").Value = -62302
$ws.Range("H132").Value = 55985.793
$ws.Range("I132").Value = 62636.938
$ws.Range("J132").Value = 3885.1667
$ws.Range("K132").Value = 187910.814
$ws.Range("L132").Value = 11655.5001
$ws.Range("M132").Value = -185380.814
$ws.Range("N132").Value = -16715.5001
$ws.Range("H133").Value = 86000
$ws.Range("J133").Value = 86000
$ws.Range("L133").Value = 86000
$ws.Range("N133").Value = -96120
